$d = $word.ActiveDocument

# The bibliography paragraph currently holds one long run of text.
# We need to insert manual line breaks (which map to <w:br/>) after:
#   "Bibliografia:"
#   "...no tópico de mentoria: "
#   "...Bullion Books, 2001. "
#   "...p. 19-27. "
#   "...Engineering Education, 2004. "
#   "...v. 21, p. 14-23,2019. "
# leaving the final reference untouched at the end.

$find = "Bibliografia:A bibliografia será recomendada pelos docentes responsáveis e obtida na busca realizada pelos próprios alunos no início dos projetos. Seguem referências no tópico de mentoria: [1] Peddy, S. The art of mentoring – Lead, follow and get out of the way. Houston: Bullion Books, 2001. [2] Zachary, L. J. The Mentor’s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promoção do bem-estar em estudantes do ensino superior. In: Programa de Monitorização e Tutorado: oito anos a promover a integração e o sucesso académico no IST. Lisboa: IST Press, 2011. p. 19-27. [3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004. [4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019. [5] Diretrizes Curriculares Nacionais para os cursos de graduação em Engenharia. Ministério da Educação. CNE/CES, 2019."

$replace = "Bibliografia:^lA bibliografia será recomendada pelos docentes responsáveis e obtida na busca realizada pelos próprios alunos no início dos projetos. Seguem referências no tópico de mentoria: ^l[1] Peddy, S. The art of mentoring – Lead, follow and get out of the way. Houston: Bullion Books, 2001. ^l[2] Zachary, L. J. The Mentor’s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promoção do bem-estar em estudantes do ensino superior. In: Programa de Monitorização e Tutorado: oito anos a promover a integração e o sucesso académico no IST. Lisboa: IST Press, 2011. p. 19-27. ^l[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004. ^l[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019. ^l[5] Diretrizes Curriculares Nacionais para os cursos de graduação em Engenharia. Ministério da Educação. CNE/CES, 2019."

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace, 2)
